$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D2").Value = "canonical SMILES"

# Column D mirrors the "canonical isomeric SMILES" column (C) for each
# microstate row, containing the (non-isomeric) canonical SMILES.
for ($r = 3; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value2
}

# Match the new column's width as recorded in the workbook XML.
# The COM layer stores ColumnWidth + 5/6 as the OOXML "width" attribute,
# so back it out here to land on the target stored width of 37.
$ws.Columns.Item(4).ColumnWidth = 37 - (5/6)
